$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the account-category labels in column E (rows 2-8) - title block
$ws.Range("E2:E8").ClearContents()

# Clear the "default values for account category" tags in column H (rows 10-67)
# These were lookup/category labels next to each trial balance line; the
# commit removes them (and their now-orphaned entries in sharedStrings).
$ws.Range("H10:H67").ClearContents()

# H69 held a check formula (=G69-F69) that is also cleared
$ws.Range("H69").ClearContents()

# Restore the view: scroll down a row and select H11:H70 (the column that
# used to hold the category tags), matching the author's last selection.
$ws.Range("H11:H70").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
